$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Write existing (reused) strings into their new positions first -------
# These reuse shared-string indices 0-9 already present in the workbook.
$ws.Range("B4").Value  = "показать отчет по депозиту"
$ws.Range("H4").Value  = "DepositExtractor"
$ws.Range("B5").Value  = "DepositViewModel"
$ws.Range("H5").Value  = "Находит все операции по данному счету"
$ws.Range("H6").Value  = "и составляет таблицу ежедневных остатков"
$ws.Range("H7").Value  = "и общие суммы взносов, процентов, расходов"
$ws.Range("B10").Value = "показать сводную форму по всем депозитам"
$ws.Range("B11").Value = "DepositsViewModel"
$ws.Range("B17").Value = "показать ожидаемые доходы от депозитов"
$ws.Range("B18").Value = "MonthAnalysisViewModel"

# --- Write brand-new strings, in the order they first appear so the ------
# shared-string table indices come out 10,11,12,... matching the target.
$ws.Range("B6").Value  = "нужна статистика и "
$ws.Range("B7").Value  = "нужен прогноз по месяцу и до конца"
$ws.Range("B12").Value = "нужна статистика и"
$ws.Range("B13").Value = "нужно определение какие %%"
$ws.Range("B14").Value = "относятся к какому году"
$ws.Range("B20").Value = "нужен прогноз по месяцу"
$ws.Range("H3").Value  = "статистика"
$ws.Range("O5").Value  = "отчеты"
$ws.Range("O6").Value  = "DepositReporter"
$ws.Range("O8").Value  = "DepositExcelReporter"
$ws.Range("O7").Value  = "составляет List<String> для отчета"
$ws.Range("O9").Value  = "составляет файл экселя"
$ws.Range("H10").Value = "расчет"
$ws.Range("H11").Value = "DepositCalculator"

# re-use of the "нужна статистика и" string (index 12) further down
$ws.Range("B19").Value = "нужна статистика и"

# --- Clear the old (pre-refactor) cell positions that are no longer used --
$ws.Range("C4").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("C7").ClearContents()
$ws.Range("C8").ClearContents()
$ws.Range("C10").ClearContents()
$ws.Range("C11").ClearContents()
$ws.Range("I4").ClearContents()
$ws.Range("I5").ClearContents()
$ws.Range("I6").ClearContents()
$ws.Range("I7").ClearContents()

# --- Bold the "class name" cells ------------------------------------------
$ws.Range("H4").Font.Bold = $true
$ws.Range("B5").Font.Bold = $true
$ws.Range("O6").Font.Bold = $true
$ws.Range("O8").Font.Bold = $true
$ws.Range("B11").Font.Bold = $true
$ws.Range("H11").Font.Bold = $true
$ws.Range("B18").Font.Bold = $true

# --- Page setup -------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Final selection ---------------------------------------------------------
$ws.Range("H11").Select() | Out-Null
